$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.012616368629892
$ws.Range("D2").Value = 1.014951748384149
$ws.Range("E2").Value = 1.014561587700151
$ws.Range("F2").Value = 1.024226717424432
$ws.Range("I2").Value = 1.025614726846036
$ws.Range("J2").Value = 1.017858058543444
$ws.Range("K2").Value = 1.017808361925375
$ws.Range("L2").Value = 1.017419369709943
$ws.Range("M2").Value = 1.027055839190686
$ws.Range("N2").Value = 1.009848667107345
$ws.Range("B3").Value = 1.019999999999999
$ws.Range("C3").Value = 1.014336710421654
$ws.Range("D3").Value = 1.0164892599011
$ws.Range("E3").Value = 1.01604693357694
$ws.Range("F3").Value = 1.025940885736105
$ws.Range("I3").Value = 1.02562406900283
$ws.Range("J3").Value = 1.019208087110867
$ws.Range("K3").Value = 1.019148899339987
$ws.Range("L3").Value = 1.01870779971146
$ws.Range("M3").Value = 1.02857458687054
$ws.Range("N3").Value = 1.010317550314682
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.01544737782087
$ws.Range("D4").Value = 1.017482030988258
$ws.Range("E4").Value = 1.017006017568406
$ws.Range("F4").Value = 1.027045749846167
$ws.Range("I4").Value = 1.025627403614609
$ws.Range("J4").Value = 1.020078946968908
$ws.Range("K4").Value = 1.02001374577063
$ws.Range("L4").Value = 1.019538984545457
$ws.Range("M4").Value = 1.029552571846218
$ws.Range("N4").Value = 1.010619419461055
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.015913716453364
$ws.Range("D5").Value = 1.017898900822223
$ws.Range("E5").Value = 1.017408740943059
$ws.Range("F5").Value = 1.027509214354519
$ws.Range("I5").Value = 1.025628156525223
$ws.Range("J5").Value = 1.02044442157123
$ws.Range("K5").Value = 1.020376722746439
$ws.Range("L5").Value = 1.019887822947322
$ws.Range("M5").Value = 1.029962592260365
$ws.Range("N5").Value = 1.010745962822157
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.015991982825428
$ws.Range("D6").Value = 1.01796886664849
$ws.Range("E6").Value = 1.017476332414375
$ws.Range("F6").Value = 1.027586972631757
$ws.Range("I6").Value = 1.02562824488384
$ws.Range("J6").Value = 1.020505749459904
$ws.Range("K6").Value = 1.020437633050004
$ws.Range("L6").Value = 1.019946360056056
$ws.Range("M6").Value = 1.030031370893023
$ws.Range("N6").Value = 1.01076718886517
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.015453611342787
$ws.Range("D7").Value = 1.017487603136696
$ws.Range("E7").Value = 1.017011400630667
$ws.Range("F7").Value = 1.027051946672684
$ws.Range("I7").Value = 1.025627416225004
$ws.Range("J7").Value = 1.020083832937655
$ws.Range("K7").Value = 1.020018598246746
$ws.Range("L7").Value = 1.01954364805179
$ws.Range("M7").Value = 1.029558054961128
$ws.Range("N7").Value = 1.010621111758523
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.013198294515499
$ws.Range("D8").Value = 1.015471800028863
$ws.Range("E8").Value = 1.015063995383784
$ws.Range("F8").Value = 1.024806930588285
$ws.Range("I8").Value = 1.025618445230591
$ws.Range("J8").Value = 1.018314872639466
$ws.Range("K8").Value = 1.018261941808386
$ws.Range("L8").Value = 1.017855327880315
$ws.Range("M8").Value = 1.027570096988831
$ws.Range("N8").Value = 1.010007447386833
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.009204262922068
$ws.Range("D9").Value = 1.011903059276468
$ws.Range("E9").Value = 1.011616323890894
$ws.Range("F9").Value = 1.020817260761663
$ws.Range("I9").Value = 1.025581881155255
$ws.Range("J9").Value = 1.015176561205916
$ws.Range("K9").Value = 1.015146301527313
$ws.Range("L9").Value = 1.014860551426519
$ws.Range("M9").Value = 1.024030175273478
$ws.Range("N9").Value = 1.008914206688083
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.006527283811466
$ws.Range("D10").Value = 1.00951195047099
$ws.Range("E10").Value = 1.009306329827307
$ws.Range("F10").Value = 1.018134002215356
$ws.Range("I10").Value = 1.025543549880452
$ws.Range("J10").Value = 1.013069402028214
$ws.Range("K10").Value = 1.013054936227039
$ws.Range("L10").Value = 1.012850093721056
$ws.Range("M10").Value = 1.021644672363744
$ws.Range("N10").Value = 1.008177147931276
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.005364524799332
$ws.Range("D11").Value = 1.008473569286867
$ws.Range("E11").Value = 1.008303175684863
$ws.Range("F11").Value = 1.016966362675369
$ws.Range("I11").Value = 1.02552364101932
$ws.Range("J11").Value = 1.012153272680669
$ws.Range("K11").Value = 1.012145810182201
$ws.Range("L11").Value = 1.011976087258644
$ws.Range("M11").Value = 1.020605491634995
$ws.Range("N11").Value = 1.007855986833456
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.004932064063981
$ws.Range("D12").Value = 1.00808740046253
$ws.Range("E12").Value = 1.007930108052128
$ws.Range("F12").Value = 1.016531765764861
$ws.Range("I12").Value = 1.02551574831709
$ws.Range("J12").Value = 1.011812409358825
$ws.Range("K12").Value = 1.011807573179005
$ws.Range("L12").Value = 1.011650908768209
$ws.Range("M12").Value = 1.020218540931337
$ws.Range("N12").Value = 1.007736386633383
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.005024854004098
$ws.Range("D13").Value = 1.00817025639555
$ws.Range("E13").Value = 1.008010152986128
$ws.Range("F13").Value = 1.016625028530645
$ws.Range("I13").Value = 1.025517463852392
$ws.Range("J13").Value = 1.01188555181133
$ws.Range("K13").Value = 1.011880151142411
$ws.Range("L13").Value = 1.011720685013948
$ws.Range("M13").Value = 1.020301586524168
$ws.Range("N13").Value = 1.007762055236113
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.005328788994164
$ws.Range("D14").Value = 1.008441658081015
$ws.Range("E14").Value = 1.008272347090997
$ws.Range("F14").Value = 1.016930456881038
$ws.Range("I14").Value = 1.025522998759562
$ws.Range("J14").Value = 1.012125108563521
$ws.Range("K14").Value = 1.012117862646188
$ws.Range("L14").Value = 1.011949218863701
$ws.Range("M14").Value = 1.020573525691029
$ws.Range("N14").Value = 1.007846106927079
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.005515978589962
$ws.Range("D15").Value = 1.008608815180387
$ws.Range("E15").Value = 1.008433833277086
$ws.Range("F15").Value = 1.017118523782957
$ws.Range("I15").Value = 1.025526343046663
$ws.Range("J15").Value = 1.012272631159239
$ws.Range("K15").Value = 1.012264251641822
$ws.Range("L15").Value = 1.012089955015729
$ws.Range("M15").Value = 1.020740949694357
$ws.Range("N15").Value = 1.007897853169119
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.006604374720592
$ws.Range("D16").Value = 1.009580799621898
$ws.Range("E16").Value = 1.009372843322618
$ws.Range("F16").Value = 1.018211371609603
$ws.Range("I16").Value = 1.025544801387166
$ws.Range("J16").Value = 1.013130123053414
$ws.Range("K16").Value = 1.013115195954131
$ws.Range("L16").Value = 1.012908024541931
$ws.Range("M16").Value = 1.021713506583702
$ws.Range("N16").Value = 1.008198419590101
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.007286117882089
$ws.Range("D17").Value = 1.010189682105837
$ws.Range("E17").Value = 1.009961070128389
$ws.Range("F17").Value = 1.01889532950999
$ws.Range("I17").Value = 1.025555493191461
$ws.Range("J17").Value = 1.013667000569308
$ws.Range("K17").Value = 1.013648010564697
$ws.Range("L17").Value = 1.013420241134484
$ws.Range("M17").Value = 1.022321884159034
$ws.Range("N17").Value = 1.008386415573326
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.007683420106069
$ws.Range("D18").Value = 1.010544543034266
$ws.Range("E18").Value = 1.010303892918888
$ws.Range("F18").Value = 1.019293714907617
$ws.Range("I18").Value = 1.025561410076097
$ws.Range("J18").Value = 1.013979794593293
$ws.Range("K18").Value = 1.01395845061309
$ws.Range("L18").Value = 1.013718674969692
$ws.Range("M18").Value = 1.022676139238386
$ws.Range("N18").Value = 1.00849587685598
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.00781883157661
$ws.Range("D19").Value = 1.010665492709685
$ws.Range("E19").Value = 1.010420739570279
$ws.Range("F19").Value = 1.019429460279187
$ws.Range("I19").Value = 1.025563373399619
$ws.Range("J19").Value = 1.014086389043117
$ws.Range("K19").Value = 1.014064245098002
$ws.Range("L19").Value = 1.01382037702934
$ws.Range("M19").Value = 1.022796829527244
$ws.Range("N19").Value = 1.0085331676296
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.007213009317406
$ws.Range("D20").Value = 1.01012438487267
$ws.Range("E20").Value = 1.009897988017736
$ws.Range("F20").Value = 1.01882200489584
$ws.Range("I20").Value = 1.025554379107118
$ws.Range("J20").Value = 1.013609435775806
$ws.Range("K20").Value = 1.013590880041246
$ws.Range("L20").Value = 1.013365319724451
$ws.Range("M20").Value = 1.022256673302155
$ws.Range("N20").Value = 1.008366265441878
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.005239303326624
$ws.Range("D21").Value = 1.008361750083004
$ws.Range("E21").Value = 1.00819515004035
$ws.Range("F21").Value = 1.016840540370336
$ws.Range("I21").Value = 1.025521382607082
$ws.Range("J21").Value = 1.012054580975586
$ws.Range("K21").Value = 1.01204787774613
$ws.Range("L21").Value = 1.011881936175436
$ws.Range("M21").Value = 1.020493472804374
$ws.Range("N21").Value = 1.007821364302786
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.003995105439615
$ws.Range("D22").Value = 1.007250797516983
$ws.Range("E22").Value = 1.007121888762575
$ws.Range("F22").Value = 1.015589594384469
$ws.Range("I22").Value = 1.025497756943018
$ws.Range("J22").Value = 1.011073665111369
$ws.Range("K22").Value = 1.01107455868811
$ws.Range("L22").Value = 1.010946180410871
$ws.Range("M22").Value = 1.01937935904964
$ws.Range("N22").Value = 1.007476986657298
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.004654992461465
$ws.Range("D23").Value = 1.00783999674686
$ws.Range("E23").Value = 1.007691097918564
$ws.Range("F23").Value = 1.016253235488904
$ws.Range("I23").Value = 1.025510554367596
$ws.Range("J23").Value = 1.011593986253805
$ws.Range("K23").Value = 1.011590838814952
$ws.Range("L23").Value = 1.011442539868393
$ws.Range("M23").Value = 1.019970500050908
$ws.Range("N23").Value = 1.007659717762466
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.007246045002887
$ws.Range("D24").Value = 1.010153890777343
$ws.Range("E24").Value = 1.009926492973108
$ws.Range("F24").Value = 1.018855138849662
$ws.Range("I24").Value = 1.02555488350111
$ws.Range("J24").Value = 1.013635447932779
$ws.Range("K24").Value = 1.013616695921315
$ws.Range("L24").Value = 1.013390137375741
$ws.Range("M24").Value = 1.022286141139235
$ws.Range("N24").Value = 1.008375371017097
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.010239266101591
$ws.Range("D25").Value = 1.012827713909171
$ws.Range("E25").Value = 1.012509611077539
$ws.Range("F25").Value = 1.021852761726391
$ws.Range("I25").Value = 1.025593793441133
$ws.Range("J25").Value = 1.01599047247932
$ws.Range("K25").Value = 1.015954232988888
$ws.Range("L25").Value = 1.015637180550932
$ws.Range("M25").Value = 1.024949778147189
$ws.Range("N25").Value = 1.009198268687294
